# AFA 2020.xlsx - Results contest 46 KKR vs KXI
# Fill in the match scores for row 55 (contest 46, "KKR vs KXI").
# The D/G/J/M/P/S columns already contain VLOOKUP/RANK formulas that
# recalculate automatically once the raw scores are entered, as do the
# SUM totals in row 70.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E55").Value = 80
$ws.Range("H55").Value = 40
$ws.Range("K55").Value = 100
$ws.Range("N55").Value = 60
$ws.Range("Q55").Value = 0
$ws.Range("T55").Value = 20

$excel.Calculate()
